# SVC-1467: add new values to lookup tables and correct spelling
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert "Wassertunnel, Wasserstollen, Druckstollen" as new row 2 ---
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "Wassertunnel, Wasserstollen, Druckstollen"
$ws.Range("B2").Value = 2023
$ws.Range("C2").Value = "http://inspire.ec.europa.eu/codelist/BuildingNatureValue/caveBuilding"
$ws.Hyperlinks.Add($ws.Range("C2"), "http://inspire.ec.europa.eu/codelist/BuildingNatureValue/caveBuilding") | Out-Null

# --- Insert "Verschlussbauwerk" as new row 9 (after Sperrwerk / before Schoepfwerk) ---
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "Verschlussbauwerk"
$ws.Range("B9").Value = 2085
$ws.Range("C9").Value = "http://inspire.ec.europa.eu/codelist/BuildingNatureValue/dam"
